# BIS-1002: removed "Internal Assignment" column from export.
# Clears the "Internal Assignment" column (O4:O8) -- header and all values --
# so the column is no longer exported, and the now-unused shared string
# "Internal Assignment" is dropped when the workbook is saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.Range("O4:O8")
$range.ClearContents()
$range.Select()
